# Append new timestamp: 2025-11-22 01:16:43 (JST) to the "取得日時" column
# on the "ランサーズ" sheet for all existing data rows (2-14).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-11-22 01:16:43"

for ($row = 2; $row -le 14; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
